$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 (pop1) with new date string and simplified "Success" sheet names
$ws.Range("B2").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"
$ws.Range("C2").Value = "Success Case Sheet.xlsx"
$ws.Range("D2").Value = "\Testdata\Templates\ImportPublications\Testing_Env\Success Case Sheet.xlsx"

# Update row 3 (pop2) with new date string and simplified "Failure" sheet names
$ws.Range("B3").Value = "UtilityOutcome - PRODFix_QOL_ECON - 9/19/2022"
$ws.Range("C3").Value = "Failure Case Sheet.xlsx"
$ws.Range("D3").Value = "\Testdata\Templates\ImportPublications\Testing_Env\Failure Case Sheet.xlsx"

# Rows 4-7 only have the date string in column B refreshed
$ws.Range("B4").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"
$ws.Range("B5").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"
$ws.Range("B6").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"
$ws.Range("B7").Value = "ICER RRMM 2022 report - ICER - 12/19/2022"

# Narrow column B to fit the new, shorter content
$ws.Columns.Item(2).ColumnWidth = 42.77

# Move the active selection from B3 to D3
$ws.Range("D3").Select()
